# Update SampleList with the new T2D samples (rows 33-61).
# - Column A (Filename) gets the new specimen ids.
# - Column H (Group (T2D or Ctrl)) gets "T2D" for every new row.
# - Formatting for the untouched columns (B..K) is copied from the last
#   existing data row (row 32) so the new rows look consistent with the
#   rest of the table.
# - The "Tableau1" table / AutoFilter range is grown to cover the new rows.
# - Final selection is left on J61, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last populated row (A32:K32) down across all
# of the new rows (A33:K61) in one shot - Excel tiles the copied block to
# fill the destination range.
$ws.Range("A32:K32").Copy()
$ws.Range("A33:K61").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$newFilenames = @(
    "C0020950", "C0020969", "C0020970", "C0020971", "C0021051", "C0021057",
    "C0023223", "C0023247", "C0023340", "C0023343", "C0023364", "C0023381",
    "C0023451", "C0023452", "C0023511", "C0023688", "C0024065", "C0024066",
    "C0024090", "C0024091", "C0024093", "C0024264", "C0024273", "C0024275",
    "C0020829", "C0020931", "C0024278", "C0024282", "C0024330"
)

$startRow = 33
for ($i = 0; $i -lt $newFilenames.Count; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $newFilenames[$i]
    $ws.Range("H$row").Value = "T2D"
}

# Grow the table / autofilter to include the newly populated rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:L61"))

# Match the final selection recorded in the source edit.
$ws.Range("J61").Select()
